$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Replace the superscript "o" in the "N°" header cell with a plain
#    degree sign "°" (no more superscript formatting).
# ------------------------------------------------------------------
$supRange = $d.Range(1, 2)
$supRange.Text = "°"
$supRange.Font.Superscript = $false

# ------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark from the end of the document to right
#    after the "N°" run. The engine refuses to add a bookmark collapsed
#    immediately in front of a hidden end-of-paragraph/cell mark, so we
#    briefly insert a placeholder character to give the insertion point
#    "room", add the bookmark there, then remove the placeholder again.
#    Re-using the "_GoBack" name makes Word move the existing bookmark
#    (the old occurrence at the end of the document disappears).
# ------------------------------------------------------------------
$placeholder = $d.Range(2, 2)
$placeholder.InsertAfter("X")

$bookmarkSpot = $d.Range(2, 2)
$d.Bookmarks.Add("_GoBack", $bookmarkSpot)

$cleanup = $d.Range(2, 3)
$cleanup.Text = ""
